# Update to tracking sheet
# Adds 5 new rows (115-119) to the "Annotations" worksheet describing
# additional people referenced by the transcriptions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Annotations")

# --- Row 115: Henry Wriothesley ---------------------------------------
$ws.Cells.Item(115, 1).Value = "Wriothesley, Henry"
$ws.Cells.Item(115, 2).Value = "Person"
$ws.Cells.Item(115, 3).Value = "psn-hwri"
$ws.Cells.Item(115, 4).Value = "../resources/annotations.xml#psn-hwri"

$note115 = "Third Earl of Southampton, Shakespeare dedicated Venus and Adonis and Rape of Lucrece to him and he has been identified as Fair Youth of Sonnets. Son of Call-Me-Wrisley from Wolf Hall..."
$cell115 = $ws.Cells.Item(115, 5)
$cell115.Value = $note115

$vaStart = $note115.IndexOf("Venus and Adonis") + 1
$vaLen = "Venus and Adonis".Length
$cell115.Characters($vaStart, $vaLen).Font.Italic = $true

$midStart = $vaStart + $vaLen
$rolStart = $note115.IndexOf("Rape of Lucrece") + 1
$midLen = $rolStart - $midStart
$cell115.Characters($midStart, $midLen).Font.Italic = $false

$rolLen = "Rape of Lucrece".Length
$cell115.Characters($rolStart, $rolLen).Font.Italic = $true

$tailStart = $rolStart + $rolLen
$tailLen = $note115.Length - $tailStart + 1
$cell115.Characters($tailStart, $tailLen).Font.Italic = $false

$ws.Rows.Item(115).RowHeight = 48

# --- Row 116: Francis Manners ------------------------------------------
$ws.Cells.Item(116, 1).Value = "Manners, Francis"
$ws.Cells.Item(116, 2).Value = "Person"
$ws.Cells.Item(116, 3).Value = "psn-fmann"
$ws.Cells.Item(116, 4).Value = "../resources/annotations.xml#psn-fmann"
$ws.Cells.Item(116, 5).Value = "Sixth Earl of Rutland, probably. Employed Shakespeare and Richard Burbage to paint his emblem."
$ws.Rows.Item(116).RowHeight = 32

# --- Row 117: Roland Whyte ----------------------------------------------
$ws.Cells.Item(117, 1).Value = "Whyte, Roland"
$ws.Cells.Item(117, 2).Value = "Person"
$ws.Cells.Item(117, 3).Value = "psn-rwhy"
$ws.Cells.Item(117, 4).Value = "../resources/annotations.xml#psn-rwhy"
$ws.Cells.Item(117, 5).Value = "https://www.jstor.org/stable/2856997"
$ws.Rows.Item(117).RowHeight = 32

# --- Row 118: Sir Robert Sidney ------------------------------------------
$ws.Cells.Item(118, 1).Value = "Sidney, Sir Robert"
$ws.Cells.Item(118, 2).Value = "Person"
$ws.Cells.Item(118, 3).Value = "psn-rsid"
$ws.Cells.Item(118, 4).Value = "../resources/annotations.xml#psn-rsid"
$ws.Cells.Item(118, 5).Value = "Younger brother to #psn-psid"
$ws.Rows.Item(118).RowHeight = 32

# --- Row 119: Mr. W. H (filled E, C, A, D to match shared-string order) --
$ws.Cells.Item(119, 5).Value = "Mysterious dedicatee of Shakespeare's Sonnets"
$ws.Cells.Item(119, 3).Value = "psn-mrwh"
$ws.Cells.Item(119, 1).Value = "Mr. W. H"
$ws.Cells.Item(119, 2).Value = "Person"
$ws.Cells.Item(119, 4).Value = "../resources/annotations.xml#psn-mrwh"
$ws.Rows.Item(119).RowHeight = 32

# --- Update the view to reflect scrolling to the newly-added rows -------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 115
$ws.Range("D119").Select()
